# Add validation checks for Excel cell formatting and formula errors
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Establish the external workbook reference (xl/externalLinks/externalLink1.xml)
# Excel needs a "real" external reference formula first so the external book
# gets registered; we then reuse it ([1]) for the real formula below and
# remove this scratch cell.
$ws.Range("AU1").Formula = "='[Sheet2.xlsx]Sheet2'!A1"
$ws.Range("AU1").ClearContents()

# --- Row 4: turn the plain error literals into real formulas that evaluate
# to the same errors.
$ws.Range("D4").Formula = "=VLOOKUP(A4,B:B,1,FALSE)"     # -> #N/A
$ws.Range("E4").Formula = "=[1]Sheet2!A1"                # -> #REF! (external)
$ws.Range("G4").Formula = "=F4"                           # -> 0
$ws.Range("I4").Formula = "=a+b"                           # -> #NAME?
$ws.Range("J4").Formula = "=K4+5"                          # -> #VALUE!

# --- Row 5: new cell with a text-formatted date-like string
$ws.Range("V5").NumberFormat = "@"
$ws.Range("V5").Value = "2/2/1902"

# --- Sheet view: move the active selection to G4
$ws.Range("G4").Select() | Out-Null

# --- Page setup: portrait orientation
$ws.PageSetup.Orientation = 1
